$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values for new rows 253-356 (drug/disease table extension) ---
$ws.Cells.Item(253, 1).Value = "ZYMRON 15 mg"
$ws.Cells.Item(253, 4).Value = "F329"
$ws.Cells.Item(254, 1).Value = "VITAMIN B6 50mg"
$ws.Cells.Item(254, 4).Value = "A154"
$ws.Cells.Item(255, 1).Value = "INH -  ISONIAZID 100 MG  TAB"
$ws.Cells.Item(255, 4).Value = "A154"
$ws.Cells.Item(256, 1).Value = "RIFAMPICIN  600  MG"
$ws.Cells.Item(256, 4).Value = "A154"
$ws.Cells.Item(257, 1).Value = "BRUFEN 400 MG"
$ws.Cells.Item(257, 4).Value = "M4796"
$ws.Cells.Item(258, 1).Value = "MYDOCALM  TAB."
$ws.Cells.Item(258, 4).Value = "M4796"
$ws.Cells.Item(259, 1).Value = "NEURONTIN 300 MG (L)"
$ws.Cells.Item(259, 4).Value = "M4796"
$ws.Cells.Item(260, 1).Value = "VITAMIN B COMPLEX"
$ws.Cells.Item(260, 4).Value = "M4796"
$ws.Cells.Item(261, 1).Value = "ZYMRON 15 mg"
$ws.Cells.Item(261, 4).Value = "F329"
$ws.Cells.Item(262, 1).Value = "SODAMINT"
$ws.Cells.Item(262, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(263, 1).Value = "CHALK CAP 350 MG"
$ws.Cells.Item(263, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(264, 1).Value = "MINOXIDIL  5 MG"
$ws.Cells.Item(264, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(265, 1).Value = "FBC TAB (FERROPRO)"
$ws.Cells.Item(265, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(266, 1).Value = "FOLIC ACID 5 MG"
$ws.Cells.Item(266, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(267, 1).Value = "CARVEDILOL 6.25 TAB"
$ws.Cells.Item(267, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(268, 1).Value = "LOSEC 20 MG"
$ws.Cells.Item(268, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(269, 1).Value = "LOPID 600 MG"
$ws.Cells.Item(269, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(270, 1).Value = "ONE- ALPHA  0.25 MCG."
$ws.Cells.Item(270, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(271, 1).Value = "VITAMIN C 500 MG 2 ML INJ."
$ws.Cells.Item(271, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(272, 1).Value = "OXYGEN  CANNULA"
$ws.Cells.Item(272, 4).Value = "D569,E782,I120,N185"
$ws.Cells.Item(273, 1).Value = "AMOXYCILLIN 500 MG"
$ws.Cells.Item(273, 4).Value = "J209"
$ws.Cells.Item(274, 1).Value = "BRICANYL 2.5 MG"
$ws.Cells.Item(274, 4).Value = "J209"
$ws.Cells.Item(275, 1).Value = "BRUFEN 400 MG"
$ws.Cells.Item(275, 4).Value = "J209"
$ws.Cells.Item(276, 1).Value = "DEXTRO 15 MG  "
$ws.Cells.Item(276, 4).Value = "J209"
$ws.Cells.Item(277, 1).Value = "PARACETAMOL 500  MG"
$ws.Cells.Item(277, 4).Value = "J209"
$ws.Cells.Item(278, 1).Value = "CAFERGOT (L)"
$ws.Cells.Item(278, 4).Value = "G439"
$ws.Cells.Item(279, 1).Value = "PARACETAMOL 500  MG"
$ws.Cells.Item(279, 4).Value = "G439"
$ws.Cells.Item(280, 1).Value = "SIBELIUM  5  MG"
$ws.Cells.Item(280, 4).Value = "G439"
$ws.Cells.Item(281, 1).Value = "OREDA"
$ws.Cells.Item(281, 4).Value = "A090"
$ws.Cells.Item(282, 1).Value = "NORFLOXACIN 400 MG"
$ws.Cells.Item(282, 4).Value = "A090"
$ws.Cells.Item(283, 1).Value = "BUSCOPAN 10 MG"
$ws.Cells.Item(283, 4).Value = "A090"
$ws.Cells.Item(284, 1).Value = "MOTILIUM 10 MG TAB"
$ws.Cells.Item(284, 4).Value = "A090"
$ws.Cells.Item(285, 1).Value = "kenalog 1 g ( oral base )"
$ws.Cells.Item(285, 4).Value = "K1379"
$ws.Cells.Item(286, 1).Value = "STUGERON 25 MG"
$ws.Cells.Item(286, 4).Value = "H813"
$ws.Cells.Item(287, 1).Value = "TAMIFLU 75 MG (L)"
$ws.Cells.Item(287, 4).Value = "J101"
$ws.Cells.Item(288, 1).Value = "BRUFEN 400 MG"
$ws.Cells.Item(288, 4).Value = "B349"
$ws.Cells.Item(289, 1).Value = "DOXYCYCLINE 100 MG"
$ws.Cells.Item(289, 4).Value = "B349"
$ws.Cells.Item(290, 1).Value = "ESSENTIALE   ( L )"
$ws.Cells.Item(290, 4).Value = "B349"
$ws.Cells.Item(291, 1).Value = "LOSEC 20 MG"
$ws.Cells.Item(291, 4).Value = "B349"
$ws.Cells.Item(292, 1).Value = "MYDOCALM  TAB."
$ws.Cells.Item(292, 4).Value = "B349"
$ws.Cells.Item(293, 1).Value = "PARACETAMOL 500  MG"
$ws.Cells.Item(293, 4).Value = "B349"
$ws.Cells.Item(294, 1).Value = "VERORAB 0.5 ML (เข็มที่4)"
$ws.Cells.Item(294, 4).Value = "Z242"
$ws.Cells.Item(295, 1).Value = "BERODUAL INH. (X)"
$ws.Cells.Item(295, 4).Value = "J439"
$ws.Cells.Item(296, 1).Value = "DEXAMETHASONE 4 MG/ 1ML INJ."
$ws.Cells.Item(296, 4).Value = "M653"
$ws.Cells.Item(297, 1).Value = "NAPROXEN 250 MG"
$ws.Cells.Item(297, 4).Value = "M653"
$ws.Cells.Item(298, 1).Value = "NORGESIC    (L)"
$ws.Cells.Item(298, 4).Value = "M653"
$ws.Cells.Item(299, 1).Value = "LOSEC 20 MG"
$ws.Cells.Item(299, 4).Value = "F410"
$ws.Cells.Item(300, 1).Value = "RIVOTRIL 2 MG"
$ws.Cells.Item(300, 4).Value = "F410"
$ws.Cells.Item(301, 1).Value = "STARIN 50 mg"
$ws.Cells.Item(301, 4).Value = "F410"
$ws.Cells.Item(302, 1).Value = "TENOFOVIR 300 mg."
$ws.Cells.Item(302, 4).Value = "B181"
$ws.Cells.Item(303, 1).Value = "CYTOTEC 200MCG"
$ws.Cells.Item(303, 4).Value = "O200"
$ws.Cells.Item(304, 1).Value = "ZOVIRAX CREAM 1 G"
$ws.Cells.Item(304, 4).Value = "A600"
$ws.Cells.Item(305, 1).Value = "POLY-OPH EYE DROP"
$ws.Cells.Item(305, 4).Value = "H001"
$ws.Cells.Item(306, 1).Value = "MAXITROL EYE OINTMENT"
$ws.Cells.Item(306, 4).Value = "H001"
$ws.Cells.Item(307, 1).Value = "EYE PADS"
$ws.Cells.Item(307, 4).Value = "H001"
$ws.Cells.Item(308, 1).Value = "CHLORAMPHENICOL EYE DROP 10 ML"
$ws.Cells.Item(308, 4).Value = "H109"
$ws.Cells.Item(309, 1).Value = "DEANXIT TAB"
$ws.Cells.Item(309, 4).Value = "F432"
$ws.Cells.Item(310, 1).Value = "ZODONREL 50 mg"
$ws.Cells.Item(310, 4).Value = "F432"
$ws.Cells.Item(311, 1).Value = "DICLOXACILLIN 500 MG."
$ws.Cells.Item(311, 4).Value = "R224"
$ws.Cells.Item(312, 1).Value = "PROZAC 20 MG"
$ws.Cells.Item(312, 4).Value = "F321"
$ws.Cells.Item(313, 1).Value = "QUANTIA 25 mg"
$ws.Cells.Item(313, 4).Value = "F321"
$ws.Cells.Item(314, 1).Value = "ATIVAN 1 MG***SA6"
$ws.Cells.Item(314, 4).Value = "F321"
$ws.Cells.Item(315, 1).Value = "BUDECORT 200 MCG/PUFF INHALER"
$ws.Cells.Item(315, 4).Value = "J459"
$ws.Cells.Item(316, 1).Value = "NSS 1000 ML ( Irrigate )=ฝาเกลียว"
$ws.Cells.Item(316, 4).Value = "J459"
$ws.Cells.Item(317, 1).Value = "ZYRTEC (L)"
$ws.Cells.Item(317, 4).Value = "J459"
$ws.Cells.Item(318, 1).Value = "DURALYN CR 200 mg"
$ws.Cells.Item(318, 4).Value = "J459"
$ws.Cells.Item(319, 1).Value = "MAXIPHED  60  MG"
$ws.Cells.Item(319, 4).Value = "J459"
$ws.Cells.Item(320, 1).Value = "STARIN 50 mg"
$ws.Cells.Item(320, 4).Value = "F320"
$ws.Cells.Item(321, 1).Value = "HALDOL 2 MG"
$ws.Cells.Item(321, 4).Value = "F320"
$ws.Cells.Item(322, 1).Value = "TRYPTANAL 25 MG*** LA1"
$ws.Cells.Item(322, 4).Value = "F320"
$ws.Cells.Item(323, 1).Value = "RIVOTRIL 2 MG"
$ws.Cells.Item(323, 4).Value = "F320"
$ws.Cells.Item(324, 1).Value = "CHLORPROMAZINE 25 MG"
$ws.Cells.Item(324, 4).Value = "F320"
$ws.Cells.Item(325, 1).Value = "AIR - X (DISFLATYL)"
$ws.Cells.Item(325, 4).Value = "K802"
$ws.Cells.Item(326, 1).Value = "MAGESTO"
$ws.Cells.Item(326, 4).Value = "K802"
$ws.Cells.Item(327, 1).Value = "BUSCOPAN 10 MG"
$ws.Cells.Item(327, 4).Value = "K802"
$ws.Cells.Item(328, 1).Value = "PROVERA 5 MG"
$ws.Cells.Item(328, 4).Value = "N938"
$ws.Cells.Item(329, 1).Value = "ZYMRON 15 mg"
$ws.Cells.Item(329, 4).Value = "G470"
$ws.Cells.Item(330, 1).Value = "RIVOTRIL 2 MG"
$ws.Cells.Item(330, 4).Value = "F328"
$ws.Cells.Item(331, 1).Value = "VALOSINE SR 75 mg"
$ws.Cells.Item(331, 4).Value = "F328"
$ws.Cells.Item(332, 1).Value = "TRANXENE 5 MG"
$ws.Cells.Item(332, 4).Value = "F328"
$ws.Cells.Item(333, 1).Value = "QUANTIA 25 mg"
$ws.Cells.Item(333, 4).Value = "F328"
$ws.Cells.Item(334, 1).Value = "MYDOCALM  TAB."
$ws.Cells.Item(334, 4).Value = "M771"
$ws.Cells.Item(335, 1).Value = "CANDID  EAR  DROPS  15 ML"
$ws.Cells.Item(335, 4).Value = "H608"
$ws.Cells.Item(336, 1).Value = "CPM  4  MG TAB"
$ws.Cells.Item(336, 4).Value = "L309"
$ws.Cells.Item(337, 1).Value = "ZYRTEC (L)"
$ws.Cells.Item(337, 4).Value = "L309"
$ws.Cells.Item(338, 1).Value = "ATARAX 10 MG*** SA3"
$ws.Cells.Item(338, 4).Value = "L239"
$ws.Cells.Item(339, 1).Value = "PREDNISOLONE 5 MG  "
$ws.Cells.Item(339, 4).Value = "L239"
$ws.Cells.Item(340, 1).Value = "BETNOVATE - N CREAM/ GM"
$ws.Cells.Item(340, 4).Value = "L239"
$ws.Cells.Item(341, 1).Value = "IV CATH NO.22"
$ws.Cells.Item(341, 4).Value = "N341"
$ws.Cells.Item(342, 1).Value = "IV SET"
$ws.Cells.Item(342, 4).Value = "N341"
$ws.Cells.Item(343, 1).Value = "NSS 100  ML INJ."
$ws.Cells.Item(343, 4).Value = "N341"
$ws.Cells.Item(344, 1).Value = "ROCEPHIN 1 G INJ. ( L )"
$ws.Cells.Item(344, 4).Value = "N341"
$ws.Cells.Item(345, 1).Value = "SYRINGE 10 CC"
$ws.Cells.Item(345, 4).Value = "N341"
$ws.Cells.Item(346, 1).Value = "NEOTICA  BALM 25 G"
$ws.Cells.Item(346, 4).Value = "M626"
$ws.Cells.Item(347, 1).Value = "DIMETAPP  TAB ( L)"
$ws.Cells.Item(347, 4).Value = "J459"
$ws.Cells.Item(348, 1).Value = "PREDNISOLONE 5 MG  "
$ws.Cells.Item(348, 4).Value = "J459"
$ws.Cells.Item(349, 1).Value = "ROPECT  "
$ws.Cells.Item(349, 4).Value = "J459"
$ws.Cells.Item(350, 1).Value = "SEROFLO 125"
$ws.Cells.Item(350, 4).Value = "J459"
$ws.Cells.Item(351, 1).Value = "ATARAX 10 MG*** SA3"
$ws.Cells.Item(351, 4).Value = "B354"
$ws.Cells.Item(352, 1).Value = "CANESTEN CREAM  1 GM"
$ws.Cells.Item(352, 4).Value = "B354"
$ws.Cells.Item(353, 1).Value = "BRUFEN 400 MG"
$ws.Cells.Item(353, 4).Value = "B029"
$ws.Cells.Item(354, 1).Value = "CPM  4  MG TAB"
$ws.Cells.Item(354, 4).Value = "B029"
$ws.Cells.Item(355, 1).Value = "VOLTAREN  25  MG TAB*** LA1/SA5"
$ws.Cells.Item(355, 4).Value = "M545"
$ws.Cells.Item(356, 1).Value = "AMOXYCILLIN 500 MG "
$ws.Cells.Item(356, 4).Value = "J304"

# --- Apply wrap-text / vertical-center style (matches existing style index 1) ---
$r = $ws.Range("A254:A284")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("A288:A293")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("A295:A301")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("A303")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("A305:A327")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("A329:A333")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("A336:A354")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("D295")
$r.WrapText = $true
$r.VerticalAlignment = -4108
$r = $ws.Range("D334")
$r.WrapText = $true
$r.VerticalAlignment = -4108

# --- Row heights for wrapped multi-line rows (ht="45") ---
$ws.Rows.Item(296).RowHeight = 45
$ws.Rows.Item(308).RowHeight = 45
$ws.Range("A315:A316").RowHeight = 45

# --- Final selection / scroll state to mirror the authored workbook ---
$ws.Range("D356").Select()
